# Refresh the crypto price/volume snapshot (columns D "Price" and E "Volume(1h)")
# with the latest scraped figures. Column D values are free-form price strings
# (often using "." as a thousands separator, e.g. "27.704.28") that must stay
# text, so they're written via a leading apostrophe (like typing '27.70 into
# the formula bar) and then the cell style is put back to "Normal" so no extra
# text-format style sticks to the cell. Column E holds padded percentage
# strings and is never misread as a number, so a plain Value assignment is
# enough there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'27.704.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").Formula = "'207.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("D8").Formula = "'22.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.06%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("D12").Formula = "'1.815.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("D13").Formula = "'1.583.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("E15").Value = "  -4.52%  "
$ws.Range("D16").Formula = "'27.693.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Formula = "'63.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Formula = "'220.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.61%  "
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Formula = "'7.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.98%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("E24").Value = "  -4.09%  "
$ws.Range("D25").Formula = "'153.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Formula = "'15.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").Formula = "'1.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("E32").Value = "  -5.06%  "
$ws.Range("D33").Formula = "'1.371.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").Formula = "'0.980"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Formula = "'0.972"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").Formula = "'64.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").Formula = "'1.726.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").Formula = "'87.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  +12.46%  "
$ws.Range("D50").Formula = "'0.0971"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("E51").Value = "  -1.16%  "
